# "code backup before pull"
#
# Adds two new QA test-case rows to the getConceptModelDataByCondition
# sheet (rows 2 & 3), widens a few columns to fit the new content, and
# flips which sheet/cell is active & selected in each sheet view.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # getConceptModelDataByCondition
$ws2 = $wb.Worksheets.Item(2)   # readDesigoCCHistoryData

# --- Seed formatting for the two new data rows from an existing,
#     already-styled data row (thin border, Calibri 10 body style) ---
$ws2.Range("A2").Copy()
$ws1.Range("A2:M3").PasteSpecial(-4122)  # xlPasteFormats

# --- Populate the new rows. Columns: A=test-id, B=description, F=name ---
$ws1.Range("F2").Value = "AM100"
$ws1.Range("F3").Value = "VS121"
$ws1.Range("A2").Value = "databrain-getConceptModelDataByCondition-test-1"
$ws1.Range("A3").Value = "databrain-getConceptModelDataByCondition-test-2"
$ws1.Range("B2").Value = "good request, data retrieved"
$ws1.Range("B3").Value = "good request, data retrieved"

# --- Widen columns A, B and F on sheet1 to fit the newly added text ---
$ws1.Columns.Item(1).ColumnWidth = 40.833333333333336
$ws1.Columns.Item(2).ColumnWidth = 22.5
$ws1.Columns.Item(6).ColumnWidth = 6

# --- Update selections on both sheets, and flip the active tab/sheet
#     from readDesigoCCHistoryData back to getConceptModelDataByCondition ---
$ws2.Activate() | Out-Null
$ws2.Range("A2").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D16").Select() | Out-Null
